$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.533148169517517
$ws.Range("B1").Value = 3.599210977554321
$ws.Range("C1").Value = 3.096096277236938
$ws.Range("D1").Value = 0.5287325978279114
$ws.Range("E1").Value = 0.8420922756195068
